# Planilha-de-custos-impressao-3D.xlsx
# "Caixa para Deck Commander - Satoru Umezawa"
#
# Update the "SEM PINTURA" cost inputs for the new print job and move the
# active-cell selection to D19 (as left by the author after editing).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUSTO")

# Filamento (kg) price
$ws.Range("D3").Value = 100
# Material utilizado (g)
$ws.Range("D4").Value = 140
# Duração da impressão (h)
$ws.Range("D5").Value = 6
# Preço sugerido (manual override used to derive lucro total)
$ws.Range("D15").Value = 45

# D9, D10, D11, D14, D16 are formulas and recalculate automatically.

# Leave the selection where the author left it after the edit.
$ws.Range("D19").Select()
